# Split the last "You know you should not think of any one." paragraph
# and add a new paragraph after it containing "Wazz up bro. ".
$d = $word.ActiveDocument

# Locate the end of the sentence we need to split after.
$anchor = $d.Content
$found = $anchor.Find.Execute("You know you should not think of any one.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor sentence in the document."
}

# Collapse to the point right after the sentence (but before the hidden
# _GoBack bookmark that lives there) and break the paragraph there.
$anchor.Collapse(0)
$anchor.InsertBefore([char]13)

# $anchor now spans the freshly inserted paragraph mark; collapse again to
# get the (empty) insertion point at the very start of the new paragraph.
$anchor.Collapse(0)

# Build the new run ("Wazz up bro. ") with the same 16pt run formatting
# used throughout the document, including the lastRenderedPageBreak marker.
$runXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r>' +
    '<w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>' +
    '<w:lastRenderedPageBreak/>' +
    '<w:t xml:space="preserve">Wazz up bro. </w:t>' +
    '</w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$anchor.InsertXML($runXml)

Write-Host "Inserted the Wazz up bro paragraph."
